$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.366.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.777.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.70%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5371'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +13.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3780'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07411'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.094'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.101'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.778.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.989'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.21%  '
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06435'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9997'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.904'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.406.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.085'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.371'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.987.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E31').Value = '  +5.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1030'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +12.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.586'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.621'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02260'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05966'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.923'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.86%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2057'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.253'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6119'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.77%  '
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('E43').Value = '  +5.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.33'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.26%  '
$ws.Range('E45').Value = '  +3.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.625'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.892'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06722'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.39%  '
